# Swap the two theme colour palettes used by this deck:
#   - the Slide Master (ppt/theme/theme1.xml) currently uses the
#     "Integral" palette and should switch to the default "Office Theme"
#     palette.
#   - the Notes Master (ppt/theme/theme2.xml) currently uses the
#     "Office Theme" palette and should switch to the "Integral" palette.
#
# PowerPoint's ColorScheme object maps 1:1 onto the 12 slots of the
# underlying DrawingML <a:clrScheme> (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink), so driving it through Colors(n).RGB reproduces the
# swap described by the diff.

$p = $ppt.ActivePresentation

function ConvertTo-BGR($hex) {
    # PowerPoint's RGB long is stored as 0x00BBGGRR.
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

$officeTheme = @("000000","FFFFFF","44546A","E7E6E6","5B9BD5","ED7D31","A5A5A5","FFC000","4472C4","70AD47","0563C1","954F72")
$integralTheme = @("000000","FFFFFF","455F51","E3DED1","99CB38","63A537","E6D024","CC9700","4EB3CF","378DA6","6B9F25","B26B02")

# Slide master (backs ppt/theme/theme1.xml) -> Office Theme palette.
$masterScheme = $p.SlideMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $masterScheme.Colors($i).RGB = ConvertTo-BGR $officeTheme[$i - 1]
}

# Notes master (backs ppt/theme/theme2.xml) -> Integral palette.
$notesScheme = $p.NotesMaster.ColorScheme
for ($i = 1; $i -le 12; $i++) {
    $notesScheme.Colors($i).RGB = ConvertTo-BGR $integralTheme[$i - 1]
}
